# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) values were recomputed from the underlying
# source data and rewritten into the worksheet. Apply the new literal
# values for each affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 2
    4  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 2
    11 = 3
    12 = 2
    13 = 2
    14 = 2
    15 = 1
    16 = 0
    17 = 0
    18 = 1
    20 = 1
    21 = 0
    22 = 0
    23 = 1
    24 = 1
    25 = 2
    26 = 1
    27 = 0
    28 = 2
    29 = 2
    30 = 0
    31 = 1
    32 = 1
    33 = 1
    34 = 1
    35 = 2
    36 = 1
    37 = 0
    39 = 2
    40 = 1
    41 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
